$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 43-49 (data range shrinks from B49 to B42)
$ws.Range("A43:B49").EntireRow.Delete()

# Update check-in / check-out date values for rows 2-42
$ws.Range("A2").Value = 45206
$ws.Range("B2").Value = 45210
$ws.Range("A3").Value = 45215
$ws.Range("B3").Value = 45216
$ws.Range("A4").Value = 45221
$ws.Range("B4").Value = 45226
$ws.Range("A5").Value = 45263
$ws.Range("B5").Value = 45267
$ws.Range("A6").Value = 45271
$ws.Range("B6").Value = 45273
$ws.Range("A7").Value = 45277
$ws.Range("B7").Value = 45279
$ws.Range("A8").Value = 45283
$ws.Range("B8").Value = 45285
$ws.Range("A9").Value = 45297
$ws.Range("B9").Value = 45301
$ws.Range("A10").Value = 45305
$ws.Range("B10").Value = 45306
$ws.Range("A11").Value = 45312
$ws.Range("B11").Value = 45313
$ws.Range("A12").Value = 45319
$ws.Range("B12").Value = 45320
$ws.Range("A13").Value = 45320
$ws.Range("B13").Value = 45325
$ws.Range("A14").Value = 45330
$ws.Range("B14").Value = 45333
$ws.Range("A15").Value = 45334
$ws.Range("B15").Value = 45338
$ws.Range("A16").Value = 45343
$ws.Range("B16").Value = 45347
$ws.Range("A17").Value = 45348
$ws.Range("B17").Value = 45350
$ws.Range("A18").Value = 45354
$ws.Range("B18").Value = 45360
$ws.Range("A19").Value = 45366
$ws.Range("B19").Value = 45368
$ws.Range("A20").Value = 45369
$ws.Range("B20").Value = 45375
$ws.Range("A21").Value = 45379
$ws.Range("B21").Value = 45383
$ws.Range("A22").Value = 45387
$ws.Range("B22").Value = 45392
$ws.Range("A23").Value = 45396
$ws.Range("B23").Value = 45399
$ws.Range("A24").Value = 45403
$ws.Range("B24").Value = 45404
$ws.Range("A25").Value = 45410
$ws.Range("B25").Value = 45411
$ws.Range("A26").Value = 45417
$ws.Range("B26").Value = 45418
$ws.Range("A27").Value = 45424
$ws.Range("B27").Value = 45425
$ws.Range("A28").Value = 45431
$ws.Range("B28").Value = 45432
$ws.Range("A29").Value = 45438
$ws.Range("B29").Value = 45439
$ws.Range("A30").Value = 45445
$ws.Range("B30").Value = 45446
$ws.Range("A31").Value = 45452
$ws.Range("B31").Value = 45453
$ws.Range("A32").Value = 45459
$ws.Range("B32").Value = 45460
$ws.Range("A33").Value = 45471
$ws.Range("B33").Value = 45473
$ws.Range("A34").Value = 45474
$ws.Range("B34").Value = 45480
$ws.Range("A35").Value = 45481
$ws.Range("B35").Value = 45485
$ws.Range("A36").Value = 45500
$ws.Range("B36").Value = 45502
$ws.Range("A37").Value = 45516
$ws.Range("B37").Value = 45520
$ws.Range("A38").Value = 45534
$ws.Range("B38").Value = 45536
$ws.Range("A39").Value = 45537
$ws.Range("B39").Value = 45543
$ws.Range("A40").Value = 45544
$ws.Range("B40").Value = 45550
$ws.Range("A41").Value = 45551
$ws.Range("B41").Value = 45557
$ws.Range("A42").Value = 45558
$ws.Range("B42").Value = 45562
